$d = $word.ActiveDocument

# Locate the paragraph that holds "Hoàn thành chức năng thêm / xóa / nhập
# lại /hiển thị chi tiết hóa đơn" (the last line under "Đã làm:" for the
# 12/11/2020 note) and insert a brand-new paragraph right after it with a
# tab followed by "Bắt lỗi", matching the surrounding formatting.
$target = "Hoàn thành chức năng thêm / xóa / nhập lại /hiển thị chi tiết hóa đơn"

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$target*") {
        $anchor = $p
    }
}

if ($anchor -ne $null) {
    # Splitting the paragraph mark here creates a new paragraph that
    # inherits the same paragraph/run formatting as $anchor.
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()

    $insertionPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
    $insertionPoint.InsertAfter("`tBắt lỗi")
}
